# SORMAS Document Template Guide — add the "Templates for Travel Entries"
# workflow section, mirroring the existing "Templates for Cases / Contacts /
# Event Participants" sections (documents and document templates for
# TravelEntries #5845).

$d = $word.ActiveDocument

# The new section belongs right after the end of the "Templates for Event
# Participants" block, i.e. after its last bullet:
#   "pathogenTest: a pathogenTest for the chosen sample (selectable)"
# That exact sentence also ends the Cases and Contacts blocks earlier in the
# document, so walk every paragraph and keep the *last* match as the anchor.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText -like "pathogenTest*chosen sample (selectable)*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'pathogenTest ... (selectable)' anchor paragraph"
}

# Create a fresh, empty paragraph right after the anchor; this is where the
# whole new block (one shaded spacer paragraph + the four new bullet
# paragraphs) will be streamed in as raw OOXML.
$anchor = $d.Paragraphs.Item($anchorIndex)
$anchor.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($anchorIndex + 1)

$travelEntriesXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Templates for Travel Entries</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>docx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>):</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Available root entities:</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>travelEntry</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>travelEntry</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> data</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>person</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: the travelEntry person</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>user</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: the current user</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newPara.Range.InsertXML($travelEntriesXml)

Write-Output "Inserted 'Templates for Travel Entries' block after paragraph $anchorIndex"
